$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Properties")
$ws.Rows.Item(19).Insert()
$ws.Range("B20:S20").Copy()
$ws.Range("B19:S19").PasteSpecial(-4122)
